# ---------------------------------------------------------------------------
# Rebuilds TestData.xlsx as a Selenium/TestNG style "test data" workbook:
#   HomePage, LoginPage, LoginPageDataProvider, RegisterPage
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the existing sheets and add the fourth one
# ------------------------------------------------------------------
$wsHome     = $wb.Worksheets.Item(1)
$wsLogin    = $wb.Worksheets.Item(2)
$wsLoginDP  = $wb.Worksheets.Item(3)

$wsHome.Name    = "HomePage"
$wsLogin.Name   = "LoginPage"
$wsLoginDP.Name = "LoginPageDataProvider"

$wsRegister = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsRegister.Name = "RegisterPage"

# Color helper value (Font.Color uses BGR ordering, not RGB) -> #202124
$colGray = 2367776

$xlFormats = -4122   # xlPasteFormats
$xlLeft    = -4131
$xlTop     = -4160

# ===================================================================
# HomePage
# ===================================================================
$wsHome.Range("A1").Value = "Title"
$wsHome.Range("B1").Value = "SubMessage"
$wsHome.Range("A2").Value = "Demo Web Shop"
$wsHome.Range("B2").Value = "Thank you for signing up! A verification email has been sent. We appreciate your interest."

# A1 -> Calibri 10, wrap text ; then propagate identical style to B1
$wsHome.Range("A1").Font.Size = 10
$wsHome.Range("A1").WrapText = $true
$wsHome.Range("A1").Copy()
$wsHome.Range("B1").PasteSpecial($xlFormats)
$excel.CutCopyMode = $false

# A2 -> Consolas 10, dark gray, left aligned, wrap text
$wsHome.Range("A2").Font.Name = "Consolas"
$wsHome.Range("A2").Font.Size = 10
$wsHome.Range("A2").Font.Color = $colGray
$wsHome.Range("A2").Font.Family = 3
$wsHome.Range("A2").HorizontalAlignment = $xlLeft
$wsHome.Range("A2").WrapText = $true

# B2 -> same Consolas font as A2, but wrap text only (no horizontal override)
$wsHome.Range("A2").Copy()
$wsHome.Range("B2").PasteSpecial($xlFormats)
$excel.CutCopyMode = $false
$wsHome.Range("B2").HorizontalAlignment = -4142   # xlGeneral -> drop the inherited "left"

$wsHome.Rows(2).RowHeight = 56.5
$wsHome.Columns("A:B").ColumnWidth = 26.63

$wsHome.PageSetup.Orientation = 1
$wsHome.Range("A2").Select()

# ===================================================================
# LoginPage
# ===================================================================
$wsLogin.Range("A1").Value = "Title"
$wsLogin.Range("B1").Value = "ErrorMessage"
$wsLogin.Range("C1").Value = "email"
$wsLogin.Range("D1").Value = "password"

$wsLogin.Range("A2").Value = "Demo Web Shop. Login"
$wsLogin.Range("B2").Value = "Login was unsuccessful. Please correct the errors and try again."
$wsLogin.Range("C2").Value = "anuj@yopmail.com"
$wsLogin.Range("D2").Value = "test@123"

# B1 -> default font, wrap text ; propagate identical style to C2
$wsLogin.Range("B1").WrapText = $true
$wsLogin.Range("B1").Copy()
$wsLogin.Range("C2").PasteSpecial($xlFormats)
$excel.CutCopyMode = $false

# A2 -> Consolas 11, dark gray, vertical top, wrap text
$wsLogin.Range("A2").Font.Name = "Consolas"
$wsLogin.Range("A2").Font.Color = $colGray
$wsLogin.Range("A2").Font.Family = 3
$wsLogin.Range("A2").VerticalAlignment = $xlTop
$wsLogin.Range("A2").WrapText = $true

# B2 -> default font, left + top aligned, wrap text
$wsLogin.Range("B2").HorizontalAlignment = $xlLeft
$wsLogin.Range("B2").VerticalAlignment = $xlTop
$wsLogin.Range("B2").WrapText = $true

$wsLogin.Columns("A").ColumnWidth = 13.36
$wsLogin.Columns("B").ColumnWidth = 20
$wsLogin.Columns("C").ColumnWidth = 19.63
$wsLogin.Rows(2).RowHeight = 32

$wsLogin.PageSetup.Orientation = 1
$wsLogin.Range("D2").Select()

# ===================================================================
# LoginPageDataProvider
# ===================================================================
$wsLoginDP.Range("B1").Value = "test@123"
$wsLoginDP.Range("A2").Value = "anuj@yopmail.com"

# A1 -> hyperlink style + vertical top + wrap, real hyperlink
$wsLoginDP.Range("A1").Value = "anuj66@yopmail.com"
$wsLoginDP.Range("A1").VerticalAlignment = $xlTop
$wsLoginDP.Range("A1").WrapText = $true
$wsLoginDP.Hyperlinks.Add($wsLoginDP.Range("A1"), "mailto:anuj66@yopmail.com") | Out-Null

# B2 -> hyperlink style, real hyperlink
$wsLoginDP.Range("B2").Value = "testing@123"
$wsLoginDP.Hyperlinks.Add($wsLoginDP.Range("B2"), "mailto:testing@123") | Out-Null

# A3 -> hyperlink style only (no live link), reuses B2's style
$wsLoginDP.Range("A3").Value = "anuj66@yopmail.com"
$wsLoginDP.Range("B2").Copy()
$wsLoginDP.Range("A3").PasteSpecial($xlFormats)
$excel.CutCopyMode = $false

# B3 -> plain text, matching value of B2/A3's data
$wsLoginDP.Range("B3").Value = "testing@123"

$wsLoginDP.Rows(1).RowHeight = 20
$wsLoginDP.Columns("A").ColumnWidth = 20.09
$wsLoginDP.Columns("B").ColumnWidth = 13.45

$wsLoginDP.Range("B24").Select()

# ===================================================================
# RegisterPage
# ===================================================================
$wsRegister.Range("A1").Value = "Gender"
$wsRegister.Range("A2").Value = "Male"
$wsRegister.Range("B1").Select()

# ------------------------------------------------------------------
# Final workbook-level view state
# ------------------------------------------------------------------
$wsLogin.Activate()

$w = $wb.Windows.Item(1)
$w.Left = 240
$w.Top = 110
$w.Width = 14810
$w.Height = 8010
